# Updates the cryptos list (refreshed Price / Volume(1h) figures, plus a
# couple of coin re-orderings and one coin swap) to match the latest
# GitHub Actions data pull.
#
# Note: several "Price" strings look like plain numbers (e.g. "298.30").
# Excel will silently coerce such text into a numeric value (dropping
# the trailing zero / decimal formatting) unless the cell is explicitly
# marked as Text beforehand, so NumberFormat = "@" is applied first for
# those specific cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.819.29'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '2.289.98'
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '298.30'
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.24'
$ws.Range("E6").Value = '  -3.29%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -2.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.54'
$ws.Range("E10").Value = '  -1.29%  '
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.57'
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.73'
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("D15").Value = '2.646.73'
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("D16").Value = '2.290.78'
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.771'
$ws.Range("E17").Value = '  -2.74%  '
$ws.Range("D18").Value = '42.733.86'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.50'
$ws.Range("E19").Value = '  -5.26%  '
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("E21").Value = '  -2.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.65'
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.59'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  -1.70%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.01'
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.41'
$ws.Range("E27").Value = '  -2.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.96'
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.53'
$ws.Range("E29").Value = '  -2.12%  '
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.76'
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.75'
$ws.Range("E34").Value = '  -3.89%  '
$ws.Range("E35").Value = '  -3.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.07'
$ws.Range("E36").Value = '  -7.07%  '
$ws.Range("E37").Value = '  -1.35%  '
$ws.Range("E38").Value = '  -1.57%  '
$ws.Range("E39").Value = '  -1.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.74'
$ws.Range("E40").Value = '  -4.35%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.73'
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.109'
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("D43").Value = '2.008.53'
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("E44").Value = '  -2.84%  '
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("E46").Value = '  +3.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.96'
$ws.Range("E47").Value = '  -2.86%  '
$ws.Range("E48").Value = '  -2.88%  '
$ws.Range("D49").Value = '2.514.04'
$ws.Range("E49").Value = '  -1.25%  '
$ws.Range("E50").Value = '  -3.48%  '
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.78'
$ws.Range("E51").Value = '  -6.68%  '